# Generate Report for Handoff
#
# 185fe346-231a-4d9d-9928-1b9795caeec0.md just got handed off (zh-cn and
# de-de), so its status flips from "In Translation" to "Ready for handoff"
# and it moves down to sit with the other "Ready for handoff" rows
# (just above 8172f727...). The rows that used to sit below it
# (1c92c8bf, bb20a7ba, aee87626) each shift up one slot to fill the gap.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A3").Value = "1c92c8bf-2a31-460d-bb19-70a900962ebf.md"
$overview.Range("A5").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
$overview.Range("A6").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
$overview.Range("A7").Value = "185fe346-231a-4d9d-9928-1b9795caeec0.md"

$overview.Range("B7").Value = "Ready for handoff"
$overview.Range("C7").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "zh-cn": A=Source File Name, B=Status, C=Latest Handoff File,
# D=Latest Handoff Datetime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A3").Value = "1c92c8bf-2a31-460d-bb19-70a900962ebf.md"
$zhcn.Range("C3").Value = "1c92c8bf-2a31-460d-bb19-70a900962ebf.7ba7829a7adeaac0d8b8e21bc6ea6fdd6a3c4464.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-02 14:39:10"

$zhcn.Range("A5").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
$zhcn.Range("C5").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.61c55838a6de6b2f7f9cc0a107a30d9e3c5128db.zh-cn.xlf"
$zhcn.Range("D5").Value = "2016-03-02 14:40:51"

$zhcn.Range("A6").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
$zhcn.Range("C6").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.c0ac426ea21c5f6307cf4c8a35e0b4b903b42ada.zh-cn.xlf"
$zhcn.Range("D6").Value = "2016-03-02 14:39:10"

$zhcn.Range("A7").Value = "185fe346-231a-4d9d-9928-1b9795caeec0.md"
$zhcn.Range("B7").Value = "Ready for handoff"
$zhcn.Range("C7").Value = "185fe346-231a-4d9d-9928-1b9795caeec0.871fef4d4bd910215edf6a3c3bd98107dcc9063e.zh-cn.xlf"
$zhcn.Range("D7").Value = "2016-03-02 14:52:36"

# ---------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A3").Value = "1c92c8bf-2a31-460d-bb19-70a900962ebf.md"
$dede.Range("C3").Value = "1c92c8bf-2a31-460d-bb19-70a900962ebf.7ba7829a7adeaac0d8b8e21bc6ea6fdd6a3c4464.de-de.xlf"
$dede.Range("D3").Value = "2016-03-02 14:39:42"

$dede.Range("A5").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.md"
$dede.Range("C5").Value = "bb20a7ba-c30e-467a-8914-ca89bacca223.61c55838a6de6b2f7f9cc0a107a30d9e3c5128db.de-de.xlf"
$dede.Range("D5").Value = "2016-03-02 14:41:01"

$dede.Range("A6").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.md"
$dede.Range("C6").Value = "aee87626-176e-4c03-92f2-8d6f0c8fb21e.c0ac426ea21c5f6307cf4c8a35e0b4b903b42ada.de-de.xlf"
$dede.Range("D6").Value = "2016-03-02 14:39:42"

$dede.Range("A7").Value = "185fe346-231a-4d9d-9928-1b9795caeec0.md"
$dede.Range("B7").Value = "Ready for handoff"
$dede.Range("C7").Value = "185fe346-231a-4d9d-9928-1b9795caeec0.871fef4d4bd910215edf6a3c3bd98107dcc9063e.de-de.xlf"
$dede.Range("D7").Value = "2016-03-02 14:52:45"
